$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'55.030.66"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "'2.295.83"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'507.36"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").Value = "'129.78"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.995"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").Value = "'0.530"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "'2.318.12"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").Value = "'0.0981"
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").Value = "'5.08"
$ws.Range("E12").Value = "  +7.46%  "
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "'23.89"
$ws.Range("E14").Value = "  +4.49%  "
$ws.Range("D15").Value = "'2.705.09"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").Value = "'54.904.96"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "'2.279.65"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "'10.71"
$ws.Range("E19").Value = "  +4.17%  "
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").Value = "'6.66"
$ws.Range("E21").Value = "  +4.07%  "
$ws.Range("D22").Value = "'310.91"
$ws.Range("E22").Value = "  +2.16%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'60.30"
$ws.Range("E24").Value = "  -2.60%  "
$ws.Range("D25").Value = "'0.992"
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +2.62%  "
$ws.Range("D28").Value = "'172.96"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'6.15"
$ws.Range("E29").Value = "  +2.99%  "
$ws.Range("E30").Value = "  +2.46%  "
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("E32").Value = "  +4.51%  "
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'0.995"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("E36").Value = "  +2.58%  "
$ws.Range("D37").Value = "'0.916"
$ws.Range("E37").Value = "  -4.92%  "
$ws.Range("D38").Value = "'3.90"
$ws.Range("E38").Value = "  +3.20%  "
$ws.Range("D39").Value = "'36.75"
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("D40").Value = "'1.44"
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("D41").Value = "'0.378"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").Value = "'135.42"
$ws.Range("E42").Value = "  +7.38%  "
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("D44").Value = "'4.91"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("D45").Value = "'259.60"
$ws.Range("E45").Value = "  +7.05%  "
$ws.Range("D46").Value = "'0.0505"
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("D47").Value = "'0.0912"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("D48").Value = "'0.552"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").Value = "'0.377"
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("E51").Value = "  +0.35%  "
